$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at row 37 (pushes old rows 37-69 down to 43-75).
# These will hold new weather entries: Light Snow, Snow, Blizzard.
$ws.Rows("37:42").Insert()

# Insert 1 new row at row 75 (after the shift above), right before the
# final "tutorialDragInstruction" row (old row 69, now at row 76).
# This will hold the new Greenland location entry.
$ws.Rows("75:75").Insert()

# --- Fill in the new cells, in the same order the values were originally
# --- typed so that new shared-string entries are appended in the exact
# --- sequence the workbook's sharedStrings.xml expects.

# New Greenland location row (Key then Value).
$ws.Range("A75").Value = "locationGreenland"
$ws.Range("B75").Value = "Greenland"

# New weather keys (column A) for rows 37-40.
$ws.Range("A37").Value = "weatherLightSnow"
$ws.Range("A38").Value = "weatherLightSnowDesc"
$ws.Range("A39").Value = "weatherSnow"
$ws.Range("A40").Value = "weatherSnowDesc"

# New weather values (column B) for rows 37-40.
$ws.Range("B37").Value = "Light Snow"
$ws.Range("B38").Value = "a dash of snow"
$ws.Range("B39").Value = "Snow"
$ws.Range("B40").Value = "a bunch of snow"

# Blizzard rows (41-42) - value/key entry order matches the source edit.
$ws.Range("B42").Value = "blizzard"
$ws.Range("A41").Value = "weatherBlizzard"
$ws.Range("A42").Value = "weatherBlizzardDesc"
$ws.Range("B41").Value = "Blizzard"

# Update the view so the active cell / visible area reflect the edit
# location (row 39, scrolled so row 25 is at the top).
$ws.Range("A39").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
